# Updates cryptocurrency price and volume figures in the cryptos worksheet
# (commit: "Updated cryptos list on Thu Oct 17 03:41:00 UTC 2024 with GitHub Actions")
#
# D-column price cells are forced to Text format before assignment (and the
# style is reset to Normal right after) so that values such as "604.75" or
# "366.00" are preserved verbatim as text instead of being auto-converted to
# numbers by Excel, matching the original inline-string cell content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.591.45"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.67%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.634.80"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.65%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.75%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.15%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.548"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.97%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.632.43"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.65%  "

$ws.Range("E10").Value = "  +7.25%  "

$ws.Range("E11").Value = "  +0.46%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.22"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.08%  "

$ws.Range("E13").Value = "  -1.19%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.04"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.88%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.113.28"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.59%  "

$ws.Range("E16").Value = "  +1.09%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.626.74"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.94%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.631.26"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.56%  "

$ws.Range("E19").Value = "  -0.14%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "366.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.72%  "

$ws.Range("E21").Value = "  -4.33%  "

$ws.Range("E22").Value = "  -0.42%  "

$ws.Range("E23").Value = "  +7.50%  "

$ws.Range("E24").Value = "  -0.05%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.05"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.19%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "66.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.91%  "

$ws.Range("E27").Value = "  +0.75%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.757.12"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "583.07"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.39%  "

$ws.Range("E30").Value = "  +0.51%  "

$ws.Range("E31").Value = "  -3.34%  "

$ws.Range("E32").Value = "  -1.23%  "

$ws.Range("E33").Value = "  -0.26%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.131"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.19%  "

$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("E36").Value = "  -2.04%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.97"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.53%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "158.07"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.06%  "

$ws.Range("E41").Value = "  -3.33%  "

$ws.Range("E42").Value = "  +0.06%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.63"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.46%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.16"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.61%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.38"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.69%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "156.43"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.11%  "

$ws.Range("E48").Value = "  -3.71%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.98"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.63%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.627"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.41%  "
